$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly reviewed/edited test case counts for
# Sprint 38 - Day 6 (rows 37-39) and Day 7 (rows 43-45)
$ws.Range("C37").Value = 579
$ws.Range("C38").Value = 824
$ws.Range("C39").Value = 615

$ws.Range("C43").Value = 629
$ws.Range("C44").Value = 854
$ws.Range("C45").Value = 615

# Update the active view/selection to reflect where the user was working
$ws.Range("C45").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
